# Weekly refresh of the Alcachofa (Mapocho Venta Directa de Santiago) data:
# the 15 data rows (rows 2-16) get reshuffled to reflect the latest source
# order. Snapshot every row first (so the row-by-row rewrite below never
# reads a value that a previous iteration already overwrote), then write
# each snapshot back into its new destination row according to the mapping
# observed in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> old row number the data should come from
$map = @{
    2  = 7
    3  = 8
    4  = 3
    5  = 4
    6  = 15
    7  = 11
    8  = 2
    9  = 16
    10 = 6
    11 = 12
    12 = 9
    13 = 14
    14 = 13
    15 = 10
    16 = 5
}

# Snapshot the whole used range (columns A-R) for every data row before
# writing anything back, since several rows swap with each other.
$snapshots = @{}
for ($r = 2; $r -le 16; $r++) {
    $snapshots[$r] = $ws.Range("A$r`:R$r").Value2
}

foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $ws.Range("A$newRow`:R$newRow").Value2 = $snapshots[$oldRow]
}
